# Apply "Add data for 2022-05-05" update:
#  - Rename sheet "Through 2022-04-26" -> "Through 2022-04-27"
#  - Update header cell I1 text "2022 (through 04-26)" -> "2022 (through 04-27)"
#  - Increment I5 (April 2022 total) from 108 -> 109
#  - Increment I14 (Total 2022 total) from 543 -> 544

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "Through 2022-04-27"

# Update the shared-string header in I1
$ws.Range("I1").Value = "2022 (through 04-27)"

# Update the two numeric cells that changed
$ws.Range("I5").Value = 109
$ws.Range("I14").Value = 544
